$wb = $excel.ActiveWorkbook

# Sheet "展览" and "全部类型" both contain the same event listing data
# in columns A:I. Update the "想去人数" (F column) figures for the
# relevant rows on both sheets.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1340
    $ws.Range("F3").Value = 1892
    $ws.Range("F4").Value = 176
    $ws.Range("F6").Value = 6316
    $ws.Range("F7").Value = 183
}
